# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all sheets
# - Shrink the Status-related columns (Overview!E:F, zh-cn!C, de-de!C)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status value wherever it currently reads "Ready for handoff"
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the status columns to their new, tighter width
$overview.Columns.Item(5).ColumnWidth = 12.55
$overview.Columns.Item(6).ColumnWidth = 12.55
$zhcn.Columns.Item(3).ColumnWidth = 12.55
$dede.Columns.Item(3).ColumnWidth = 12.55
